$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "293.24"
Set-TextValue "E2" "2.34%"
Set-TextValue "D3" "30.57"
Set-TextValue "E3" "6.87%"
Set-TextValue "D4" "5.152"
Set-TextValue "E4" "1.76%"
Set-TextValue "D5" "0.07131"
Set-TextValue "E5" "7.07%"
Set-TextValue "D6" "7.535"
Set-TextValue "E6" "1.97%"
Set-TextValue "D7" "3.631"
Set-TextValue "E7" "6.54%"
Set-TextValue "E8" "2.50%"
Set-TextValue "D9" "0.9169"
Set-TextValue "E9" "-1.96%"
Set-TextValue "D10" "0.1630"
Set-TextValue "E10" "3.36%"
Set-TextValue "D11" "0.07688"
Set-TextValue "E11" "18.77%"
Set-TextValue "D12" "0.07758"
Set-TextValue "E12" "2.62%"
Set-TextValue "E13" "0.18%"
Set-TextValue "E14" "0.20%"
Set-TextValue "D15" "0.001580"
Set-TextValue "E15" "-0.34%"
Set-TextValue "D16" "0.0006556"
Set-TextValue "E16" "1.34%"
Set-TextValue "D17" "0.006509"
Set-TextValue "E17" "3.92%"
Set-TextValue "D18" "3.483"
Set-TextValue "E18" "1.09%"
Set-TextValue "D19" "2.239"
Set-TextValue "E19" "-0.48%"
Set-TextValue "D20" "0.3251"
Set-TextValue "E20" "1.04%"
Set-TextValue "D21" "0.1365"
Set-TextValue "E21" "5.19%"
Set-TextValue "D22" "3.847"
Set-TextValue "E22" "-5.88%"
Set-TextValue "D23" "0.1601"
Set-TextValue "E23" "3.17%"
Set-TextValue "D24" "0.04530"
Set-TextValue "E24" "0.71%"
Set-TextValue "E25" "2.63%"
Set-TextValue "D26" "0.004236"
Set-TextValue "E26" "2.25%"
Set-TextValue "E27" "-6.32%"
Set-TextValue "D28" "0.0001691"
Set-TextValue "E28" "4.55%"
Set-TextValue "D40" "0.04410"
Set-TextValue "E40" "4.67%"
Set-TextValue "D41" "0.007020"
Set-TextValue "E41" "4.34%"
Set-TextValue "D42" "0.1272"
Set-TextValue "E42" "1.98%"
Set-TextValue "E43" "9.53%"
Set-TextValue "D44" "0.01320"
Set-TextValue "E44" "8.50%"
Set-TextValue "D45" "0.00005872"
Set-TextValue "E45" "4.47%"
Set-TextValue "D47" "0.01300"
Set-TextValue "E47" "-0.44%"
